$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = -400
$ws.Range("E20").Value = -600
$ws.Range("F20").Value = -1800

$ws.Range("D21").Value = -3300
$ws.Range("E21").Value = -4600
$ws.Range("F21").Value = -6600

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0

$ws.Range("D27").Value = -6700

$ws.Range("D32").Value = 400
$ws.Range("E32").Value = 600
$ws.Range("F32").Value = 1800

$ws.Range("D33").Value = -6700

$ws.Range("D35").Value = -6700

$ws.Range("D81").Value = -6700

$ws.Range("D91").Value = -6600
$ws.Range("E91").Value = -3700
$ws.Range("F91").Value = -1000
